$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds date-like strings (e.g. "2026-02-20") that must stay plain
# text, not get auto-converted to Excel date serials. Force text format,
# set the values, then reset the style so no extra formatting sticks.
$dateRng = $ws.Range("F2:F5")
$dateRng.NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Senior Data Engineer"
$ws.Range("B2").Value = "Deloitte"
$ws.Range("C2").Value = "San Jose, CA, US USA"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = "RAG, Glue, Redshift, BigQuery, Synapse, Git, Snowflake, Databricks, BigQuery, Redshift"
$ws.Range("F2").Value = "2026-02-20"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=8ed48e228bf802be"

# Row 3
$ws.Range("A3").Value = "Software Engineer - Database Integrations"
$ws.Range("B3").Value = "clickhouse"
$ws.Range("C3").Value = "Remote, US USA"
$ws.Range("D3").Value = 12.2
$ws.Range("E3").Value = "RAG, BigQuery, Data Lake, Kubernetes, Snowflake, BigQuery, Kafka, MySQL, MongoDB, SQL"
$ws.Range("F3").Value = "2026-02-21"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=77453e14526da77d"

# Row 4
$ws.Range("A4").Value = "AI & Data Consultant"
$ws.Range("B4").Value = "Deloitte"
$ws.Range("C4").Value = "Chicago, IL, US USA"
$ws.Range("D4").Value = 12.2
$ws.Range("E4").Value = "Data Scientist, LangChain, RAG, CI/CD, Git, Tableau, Power BI, Python, SQL, R"
$ws.Range("F4").Value = "2026-02-20"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=a97fe7cbba5b42a5"

# Row 5
$ws.Range("A5").Value = "Sr. Business Intelligence Engineer - Digital Experiences & Capabilities"
$ws.Range("B5").Value = "Visa"
$ws.Range("C5").Value = "San Francisco, CA, US USA"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = "RAG, Git, Hadoop, Tableau, Power BI, R, Scala, Optimization, A/B Testing"
$ws.Range("F5").Value = "2026-02-20"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=3df7a7234ca3a31b"

# Restore the default (unstyled) look for the date column now that the
# values are locked in as text.
$dateRng.Style = "Normal"

# Remove row 6 entirely (job list shrank from 5 to 4 entries)
$ws.Rows.Item(6).Delete()
